$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily drop the AutoFilter so the subsequent column insert doesn't
# get confused by the existing filtered range, then rebuild it afterwards
# over the new, wider range.
$ws.AutoFilterMode = $false

# Insert a new blank column before column S ("Brand" is R, the old
# "Logical Operator" column was S) - this shifts S:AO -> T:AP and extends
# the used range / long per-row style-21 tail from AML to AMM.
$ws.Columns("S").Insert()

# Label the newly inserted column and put it in the shared-string table.
$ws.Range("S1").Value = "Sub brand"

# Recreate the AutoFilter over the new, one-column-wider range.
$ws.Range("A1:AP54").AutoFilter()

# The two _FilterDatabase defined names still point at the old AO54 range;
# repoint them at the new AP54 range.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Cinema!`$A`$1:`$AP`$54"
$wb.Names.Item("_FilterDatabase_0").RefersTo = "=Cinema!`$A`$1:`$AP`$54"

# Restore/move the active selection to the newly inserted column's data
# cell, matching the saved cursor position.
$ws.Range("S2").Select()
